$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Clear out the old "Test Scanner/Parse/Names/GUI" rows (8-11) and the
#    trailing empty rows (13, 14) so stale shared strings get pruned and the
#    row/column layout can be rebuilt cleanly.
# ---------------------------------------------------------------------------
$ws.Range("B4:N14").Clear()

# ---------------------------------------------------------------------------
# 2) Re-write column B task labels (rows 4-17) in the new order.
#    Reused strings (Write Scanner/Parse/Names/GUI Module, Maintenance) keep
#    their existing shared-string slots; brand-new strings are introduced in
#    the exact order needed to reproduce the target shared string table.
# ---------------------------------------------------------------------------

# -- legend first (introduces James / Anna / Neelay ahead of the task rows) --
$ws.Range("C19").Value = "James"
$ws.Range("D19").Value = "Anna"
$ws.Range("E19").Value = "Neelay"

$ws.Range("B4").Value = "Write Names Module"
$ws.Range("B5").Value = "Write Scanner Module"
$ws.Range("B6").Value = "Write Names Unit Tests"
$ws.Range("B7").Value = "Write Scannner Unit Tests"
$ws.Range("B8").Value = "Debug Names"
$ws.Range("B9").Value = "Debug Scanner"
$ws.Range("B10").Value = "Write Parse Module"
$ws.Range("B11").Value = "Write GUI Module"
$ws.Range("B12").Value = "Write Parse Unit Tests"
$ws.Range("B13").Value = "Write GUI Unit Tests"
$ws.Range("B14").Value = "Debug Parse"
$ws.Range("B15").Value = "Debug GUI"
$ws.Range("B16").Value = "Integrate Modules"
$ws.Range("B17").Value = "Maintenance"

$ws.Range("F19").Value = "All"

# ---------------------------------------------------------------------------
# 3) Re-add the N-column marker cells (N4:N14) that Clear() removed.
# ---------------------------------------------------------------------------
$ws.Range("N4").Value = ""
$ws.Range("N5").Value = ""
$ws.Range("N6").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("N8").Value = ""
$ws.Range("N9").Value = ""
$ws.Range("N10").Value = ""
$ws.Range("N11").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("N14").Value = ""
$ws.Range("N4:N14").Borders.LineStyle = 1
$ws.Range("N4:N14").Borders.LineStyle = 0

Write-Output "content pass done"
